# Update the patient "id" value on row 2 of Sheet1 (EEG sample data),
# matching the CA-xxxxxxxx identifiers used across the genetics/cgm/eeg
# fixtures. The final id stored in A2 becomes "CA-AUV5IQQD".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CA-AUV5IQQD"
